# Weight Data.xlsx - add a new "11/24/2023" entry at the top of the log
# and drop the now-superseded "11/03/2023" entry (merge-tables sequence
# tidy-up per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row right under the header (row 1) ---
$ws.Rows.Item(2).Insert()

# Column A holds date-like text ("11/24/2023"). A plain .Value assignment
# gets auto-parsed into a real date serial by Excel, which also forces a
# new number-format style onto the cell - neither of which matches the
# sheet (every data cell is plain shared-string text with the default
# style). Route it through a scratch formula cell + copy / paste-values
# so it lands back as literal text with no style change.
$scratch = $ws.Range("ZZ1000")
$scratch.Formula = "=""11/24/2023"""
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B2").Value = "15.3% Fat"
$ws.Range("C2").Value = "70.9 kg"

# --- Remove the stale "11/03/2023" row ---
# After the insert above, every original row shifted down by one, so the
# row that used to be 21 (11/03/2023 | 16.2% Fat | 72.7 kg) is now row 22.
$ws.Rows.Item(22).Delete()
